$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 144 ("Camote"/"Zapallo" detail
# rows). This shifts the old rows 144-151 down to 147-154 and leaves three
# fresh blank rows at 144-146 for the new weekly records.
$ws.Rows("144:146").Insert()

# --- Row 144 (new record) ---
$ws.Cells.Item(144, 1).Value = 5
$ws.Cells.Item(144, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(144, 3).Value = "Maule"
$ws.Cells.Item(144, 4).Value = 44516
$ws.Cells.Item(144, 5).Value = 7
$ws.Cells.Item(144, 6).Value = 100112045
$ws.Cells.Item(144, 7).Value = "Zapallo"
$ws.Cells.Item(144, 8).Value = "Camote"
$ws.Cells.Item(144, 9).Value = "1a (guarda)"
$ws.Cells.Item(144, 10).Value = 800
$ws.Cells.Item(144, 11).Value = 500
$ws.Cells.Item(144, 12).Value = 500
$ws.Cells.Item(144, 13).Value = 500
$ws.Cells.Item(144, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(144, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(144, 16).Value = 500
$ws.Cells.Item(144, 17).Value = 1
$ws.Cells.Item(144, 18).Value = "Hortaliza"

# --- Row 145 (new record) ---
$ws.Cells.Item(145, 1).Value = 5
$ws.Cells.Item(145, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(145, 3).Value = "Maule"
$ws.Cells.Item(145, 4).Value = 44516
$ws.Cells.Item(145, 5).Value = 7
$ws.Cells.Item(145, 6).Value = 100112045
$ws.Cells.Item(145, 7).Value = "Zapallo"
$ws.Cells.Item(145, 8).Value = "Camote"
$ws.Cells.Item(145, 9).Value = "1a nueva(o)"
$ws.Cells.Item(145, 10).Value = 800
$ws.Cells.Item(145, 11).Value = 650
$ws.Cells.Item(145, 12).Value = 650
$ws.Cells.Item(145, 13).Value = 650
$ws.Cells.Item(145, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(145, 15).Value = "Perú"
$ws.Cells.Item(145, 16).Value = 650
$ws.Cells.Item(145, 17).Value = 1
$ws.Cells.Item(145, 18).Value = "Hortaliza"

# --- Row 146 (new record) ---
$ws.Cells.Item(146, 1).Value = 5
$ws.Cells.Item(146, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(146, 3).Value = "Maule"
$ws.Cells.Item(146, 4).Value = 44516
$ws.Cells.Item(146, 5).Value = 7
$ws.Cells.Item(146, 6).Value = 100112045
$ws.Cells.Item(146, 7).Value = "Zapallo"
$ws.Cells.Item(146, 8).Value = "Paine"
$ws.Cells.Item(146, 9).Value = "1a (guarda)"
$ws.Cells.Item(146, 10).Value = 2500
$ws.Cells.Item(146, 11).Value = 80
$ws.Cells.Item(146, 12).Value = 80
$ws.Cells.Item(146, 13).Value = 80
$ws.Cells.Item(146, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(146, 15).Value = "Región del Maule"
$ws.Cells.Item(146, 16).Value = 80
$ws.Cells.Item(146, 17).Value = 1
$ws.Cells.Item(146, 18).Value = "Hortaliza"
